$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "doctor_MA" average column (AF) values for rows 4-13
$ws.Range("AF4").Value = 0.5620000000000001
$ws.Range("AF5").Value = 0.75
$ws.Range("AF6").Value = 0.643
$ws.Range("AF7").Value = 0.703
$ws.Range("AF8").Value = 0.75
$ws.Range("AF9").Value = 0.75
$ws.Range("AF10").Value = 0.75
$ws.Range("AF11").Value = 0.75
$ws.Range("AF12").Value = 1
$ws.Range("AF13").Value = 2
